$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.828.47"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.112.79"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.98"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.82"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.108.90"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("E10").Value = "  -2.35%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.05"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "3.628.54"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "66.801.35"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.17"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("D19").Value = "3.113.24"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.28"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "475.62"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.713"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.87"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.37"
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.55"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "0.0₃0938"
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.974"
$ws.Range("E37").Value = "  -2.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "46.94"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.13"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("E41").Value = "  -2.53%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.70"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.59"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.83"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  -1.76%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.10"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.33"
$ws.Range("E25").Value = "  +3.71%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "386.91"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.827.82"
$ws.Range("E45").Value = "  +2.60%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0355"
$ws.Range("E46").Value = "  -1.39%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  -8.34%  "
